# Cargo Dist Conversion Factors.xlsx — "Uploading newest EPS-US files"
#
# Switches the CDCF workbook's "About" notes and the two conversion-factor
# sheets from the EU (metric, Gtkm/Gpkm) output-unit convention to the US
# (trillion passenger-miles / trillion freight ton-miles) convention, and
# updates the two conversion formulas from unit-specific multipliers to a
# straight 10^12 (since the US output units are already in the model's
# native miles-based units, just scaled to trillions).

$wb = $excel.ActiveWorkbook

# ---- "About" sheet -------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop the old EU-only "ton-mile conversion" note (row 15) and the helper
# formula cell below it (row 17, after the row-15 delete it is row 16) —
# delete bottom row first so the row numbers above stay valid.
$about.Rows.Item(17).Delete()
$about.Rows.Item(15).Delete()

# Re-word the "desired output units" blurb for the US model.
$about.Range("A11").Value = "For the U.S. model, the desired output units are:"
$about.Range("A12").Value = "trillion passenger-miles"
$about.Range("A13").Value = "trillion freight ton-miles"

# ---- "CDCF-PMpPDOU" sheet ------------------------------------------
$pm = $wb.Worksheets.Item("CDCF-PMpPDOU")
$pm.Range("B2").Formula = "=10^12"

# ---- "CDCF-FTMpFDOU" sheet -----------------------------------------
$ftm = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ftm.Range("B2").Formula = "=10^12"
# Remove the leftover formatted-but-empty cell below the data (row 5).
$ftm.Rows.Item(5).Delete()
